$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "StatQuery" Cypher text for the FilesTab row.
# This same shared string is also used by the CasesTab (C2) and
# SamplesTab (C3) rows, so update all three cells to keep them in sync.
$newQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['Australian Shepherd']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$ws.Range("C2:C4").Value = $newQuery

# Scroll back to the top of the sheet and move the active selection to B1
# (previously the view was scrolled down with B4 selected).
$window = $excel.ActiveWindow
$window.ScrollRow = 1
$window.ScrollColumn = 1
$ws.Range("B1").Select()
